$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Resource1-230120 (Chinese)"
$ws.Range("B2").Value = "Resource8-230120 Summary"

$ws.Range("A3").Value = "Resource2-230120 (Private Group)"
$ws.Range("B3").Value = "Resource2-230120 Summary"

$ws.Range("A4").Value = "Resource3-230120 (External Members)"
$ws.Range("B4").Value = "Resource3-230120 Summary"

$ws.Range("A5").Value = "Resource4-230120 (Network)"
$ws.Range("B5").Value = "Resource4-230120 Summary"

$ws.Range("A6").Value = "Resource5-230120 (Distributed Resiliency)"
$ws.Range("B6").Value = "Resource5-230120 Summary"

$ws.Range("A7").Value = "Resource6-230120 (Internal Admins)"
$ws.Range("B7").Value = "Resource6-230120 Summary"

$ws.Range("A8").Value = "Resource7-230120 (Private Group)"
$ws.Range("B8").Value = "Resource7-230120 Summary"
